$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2058823529411765
$ws.Range("C2").Value = 0.5
$ws.Range("P2").Value = 0.08823529411764706
$ws.Range("S2").Value = 0.2058823529411765
$ws.Range("J3").Value = 0.05882352941176471
$ws.Range("P3").Value = 0.5882352941176471
$ws.Range("S3").Value = 0.3529411764705883
$ws.Range("S4").Value = 1
$ws.Range("B6").Value = 0.03846153846153846
$ws.Range("D6").Value = 0.03846153846153846
$ws.Range("F6").Value = 0.07692307692307693
$ws.Range("J6").Value = 0.1923076923076923
$ws.Range("Q6").Value = 0.1538461538461539
$ws.Range("R6").Value = 0.03846153846153846
$ws.Range("S6").Value = 0.4615384615384616
$ws.Range("B7").Value = 0.02857142857142857
$ws.Range("F7").Value = 0.08571428571428572
$ws.Range("J7").Value = 0.1142857142857143
$ws.Range("Q7").Value = 0.1428571428571428
$ws.Range("S7").Value = 0.6285714285714286
$ws.Range("B8").Value = 0.08
$ws.Range("F8").Value = 0.04
$ws.Range("J8").Value = 0.08
$ws.Range("O8").Value = 0.06666666666666667
$ws.Range("Q8").Value = 0.2533333333333334
$ws.Range("R8").Value = 0.01333333333333333
$ws.Range("S8").Value = 0.4666666666666667
$ws.Range("F9").Value = 0.05405405405405406
$ws.Range("J9").Value = 0.1621621621621622
$ws.Range("Q9").Value = 0.1081081081081081
$ws.Range("R9").Value = 0.02702702702702703
$ws.Range("S9").Value = 0.6486486486486487
$ws.Range("B10").Value = 0.09844559585492228
$ws.Range("F10").Value = 0.04663212435233161
$ws.Range("J10").Value = 0.1243523316062176
$ws.Range("Q10").Value = 0.2642487046632124
$ws.Range("R10").Value = 0.02072538860103627
$ws.Range("S10").Value = 0.4455958549222798
$ws.Range("G11").Value = 0.1111111111111111
$ws.Range("J11").Value = 0.04444444444444445
$ws.Range("K11").Value = 0.1333333333333333
$ws.Range("L11").Value = 0.7111111111111111
$ws.Range("G12").Value = 0.7575757575757576
$ws.Range("J12").Value = 0.2424242424242424
$ws.Range("G13").Value = 0.6
$ws.Range("J13").Value = 0.3
$ws.Range("S13").Value = 0.1
$ws.Range("H15").Value = 0.2
$ws.Range("I15").Value = 0.1666666666666667
$ws.Range("J15").Value = 0.2666666666666667
$ws.Range("K15").Value = 0.03333333333333333
$ws.Range("M15").Value = 0.03333333333333333
$ws.Range("S15").Value = 0.3
$ws.Range("I16").Value = 0.08333333333333333
$ws.Range("J16").Value = 0.6666666666666666
$ws.Range("K16").Value = 0.08333333333333333
$ws.Range("F17").Value = 0.01204819277108434
$ws.Range("H17").Value = 0.1807228915662651
$ws.Range("I17").Value = 0.1325301204819277
$ws.Range("J17").Value = 0.3855421686746988
$ws.Range("K17").Value = 0.108433734939759
$ws.Range("O17").Value = 0.08433734939759036
$ws.Range("S17").Value = 0.0963855421686747
$ws.Range("F18").Value = 0.1428571428571428
$ws.Range("H18").Value = 0.1428571428571428
$ws.Range("J18").Value = 0.2857142857142857
$ws.Range("K18").Value = 0.2857142857142857
$ws.Range("S18").Value = 0.1428571428571428
$ws.Range("F19").Value = 0.008438818565400843
$ws.Range("H19").Value = 0.219409282700422
$ws.Range("I19").Value = 0.08438818565400844
$ws.Range("J19").Value = 0.3713080168776371
$ws.Range("K19").Value = 0.1054852320675106
$ws.Range("M19").Value = 0.0379746835443038
$ws.Range("O19").Value = 0.05485232067510549
$ws.Range("S19").Value = 0.1181434599156118
